# Generate Report for Handoff
#
# The localization item moved out of the "handed back" state and is now
# queued for a fresh handoff, so:
#   - every "Status" cell that read "Handed back: in sync with en-US"
#     now reads "Ready for handoff"
#   - the handoff-generation timestamps that were stamped at the old
#     status change are refreshed to the new run's time
#   - the (now much shorter) status columns are narrowed to fit the new
#     text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refresh timestamps recorded at handoff-generation time
# Overview!G2 ("Latest HO Xliff Generate Date") and de-de!H2 ("Latest
# Handoff Datetime") shared the same value before the edit, and still do
# after it, so both need to be written.
$overview.Range("G2").Value = "2016-08-25 13:02:21"
$dede.Range("H2").Value = "2016-08-25 13:02:21"

# zh-cn!H2 ("Latest Handoff Datetime")
$zhcn.Range("H2").Value = "2016-08-25 13:02:16"

# --- Narrow the status columns now that "Ready for handoff" is shorter
# than "Handed back: in sync with en-US". 16.33 is the ColumnWidth
# (character units) that Excel stores as the same pixel-rounded column
# width used in the refreshed report.
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
